$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Tasks estimation"
$ws2 = $wb.Worksheets.Item(2)   # "Costs"

# --- Tasks estimation sheet edits ---

# Remove the two consecutive rows "Login via Facebook & Google" (row 14)
# and "E-mail notification" (row 15).
$ws1.Range("A14:A15").EntireRow.Delete()

# After the above deletion, "Contact administrator (notification via e-mail)"
# has shifted up from row 19 to row 17 - remove it as well.
$ws1.Range("A17").EntireRow.Delete()

# Update the re-estimated hour counts.
$ws1.Range("C8").Value = 32
$ws1.Range("C10").Value = 32

# Insert a new task row right after "Add comments (design page)" (row 14)
# for the new "Add reply comments (design page)" task.
$ws1.Range("A15").EntireRow.Insert()
$ws1.Range("B15").Value = "Add reply comments (design page)"
$ws1.Range("C15").Value = 20

# --- Active sheet / selection bookkeeping ---

# Selection on the Costs sheet moves to B3.
$ws2.Range("B3").Select()

# The Tasks estimation sheet becomes the active tab, selection at C9.
$ws1.Activate()
$ws1.Range("C9").Select()
